$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

# Remove the two leading "Properties"/"Value" placeholder rows so the
# "Phase/Run/Gmail/..." header becomes row 1 and the data shifts up.
$ws.Rows("1:1").Select() | Out-Null
$ws.Rows("1:2").Delete() | Out-Null
